$d = $word.ActiveDocument

$replacements = @(
    @{old = "760÷7=108, 4"; new = "431÷2=215, 1"},
    @{old = "683÷4=170, 3"; new = "697÷3=232, 1"},
    @{old = "163÷8=20, 3"; new = "357÷2=178, 1"},
    @{old = "664÷8=83, 0"; new = "737÷3=245, 2"},
    @{old = "196÷9=21, 7"; new = "286÷9=31, 7"},
    @{old = "980÷4=245, 0"; new = "436÷5=87, 1"},
    @{old = "164÷8=20, 4"; new = "541÷4=135, 1"},
    @{old = "873÷3=291, 0"; new = "378÷5=75, 3"},
    @{old = "746÷2=373, 0"; new = "713÷7=101, 6"},
    @{old = "390÷9=43, 3"; new = "472÷7=67, 3"},
    @{old = "570÷4=142, 2"; new = "191÷8=23, 7"},
    @{old = "959÷9=106, 5"; new = "520÷8=65, 0"},
    @{old = "927÷7=132, 3"; new = "924÷7=132, 0"},
    @{old = "358÷6=59, 4"; new = "504÷5=100, 4"},
    @{old = "828÷4=207, 0"; new = "291÷5=58, 1"},
    @{old = "741÷8=92, 5"; new = "370÷8=46, 2"},
    @{old = "432÷8=54, 0"; new = "613÷7=87, 4"},
    @{old = "534÷7=76, 2"; new = "937÷5=187, 2"},
    @{old = "820÷4=205, 0"; new = "219÷6=36, 3"},
    @{old = "542÷5=108, 2"; new = "997÷2=498, 1"},
    @{old = "910÷3=303, 1"; new = "845÷6=140, 5"},
    @{old = "486÷2=243, 0"; new = "598÷8=74, 6"},
    @{old = "872÷6=145, 2"; new = "979÷8=122, 3"},
    @{old = "710÷6=118, 2"; new = "499÷5=99, 4"},
    @{old = "545÷3=181, 2"; new = "149÷8=18, 5"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
